# Settings back to qfin22 presentation - recalibration
# Update calibration parameters on the AR, SETAR, GARCH, TARCH and AR_TARCH
# sheets of the riskDriverType/PnL factor calibration workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR")
$ws.Range("B2").Value = 0.08416508657525568
$ws.Range("B3").Value = 0.7855789005287703
$ws.Range("B4").Value = 364.6169178066191
$ws.Range("B5").Value = "[1.0, 0.1757729194359353, 0.13012079664062226, 0.13285524394835732, 0.171231486293588, 0.3752979122887914, 0.1549521660631866, 0.18932627239660701, 0.16266357456861272, 0.16422928049539764, 0.1240252489274742, 0.13463090455052687, 0.12976949264902557, 0.17904142714016563, 0.21926570175479085, 0.12165396046352103, 0.15867344409494355, 0.1275472986787722, 0.10448933163764733, 0.12222005011539219]"

# ---------------------------------------------------------------------------
# SETAR sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("SETAR")
$ws.Range("B4").Value = -6.702807285249546
$ws.Range("B5").Value = 0.6936946348852393
$ws.Range("B6").Value = 233.616705020562
$ws.Range("B7").Value = 6.731080728003657
$ws.Range("B8").Value = 0.702354320098352
$ws.Range("B9").Value = 234.2150180554948
$ws.Range("B10").Value = "[1.0, 0.23287956540147117, 0.17443777814252565, 0.16634614541679668, 0.20734471860375378, 0.26343586948573977, 0.1317104681130148, 0.14649607967987363, 0.12047166290631892, 0.15988173856564072, 0.1649064530858388, 0.1353254875036473, 0.1129377167035112, 0.13788157045406166, 0.16828732077259884, 0.10635784535691122, 0.14230833922483252, 0.1489336327778308, 0.12163180305889687, 0.08738990287158936]"

# ---------------------------------------------------------------------------
# GARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GARCH")
$ws.Range("B2").Value = 0.01710343389510583
$ws.Range("B3").Value = 3.549083294622785
$ws.Range("B4").Value = 0.04391208679432808
$ws.Range("B5").Value = 0.945458116854444
$ws.Range("B6").Value = "[0.9999999999999998, 0.010988428378977882, -0.04921181731271194, -0.06420760043620628, -0.05273345997281751, 0.29035375354031645, 0.025195080282601263, 0.007228314256517425, 0.01719153663398028, -0.006857227427753952, -0.06907271766159237, -0.06106488540700574, -0.030598653348629365, 0.023508173038901327, 0.046757296424504845, -0.014689149774101308, -0.018507088069220308, -0.04268243504341597, -0.04711210901952304, -0.035482238107404585]"

# ---------------------------------------------------------------------------
# TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TARCH")
$ws.Range("B2").Value = -0.09771786591627307
$ws.Range("B3").Value = 3.502367387825018
$ws.Range("B4").Value = 0.03520051902499573
$ws.Range("B5").Value = 0.9454089928665422
$ws.Range("B6").Value = "[0.9999999999999998, 0.010165607450609147, -0.04901309950757094, -0.06474251635158627, -0.052745794558228046, 0.2905333561546549, 0.024978388249056033, 0.007295603821316108, 0.01638566913349809, -0.006948849840473221, -0.06814324204146022, -0.06108465973747744, -0.030650648877121347, 0.023370607485687693, 0.04710307635130273, -0.014487978622320416, -0.017955918268259645, -0.041886205254117084, -0.04650252755448629, -0.035542278975759706]"
$ws.Range("B7").Value = 0.0181493234319557

# ---------------------------------------------------------------------------
# AR_TARCH sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("AR_TARCH")
$ws.Range("B2").Value = 0.01909915648711362
$ws.Range("B3").Value = 3.188892034704653
$ws.Range("B4").Value = 0.03550416484973991
$ws.Range("B5").Value = 0.948439846431152
$ws.Range("B6").Value = "[0.9999999999999998, 0.02452513744791199, -0.038075545904617855, -0.05876004774555581, -0.014199634878138497, 0.22091182599911535, 0.021481848965776852, 0.014425549944511217, -0.011919078292818026, -0.009683242255229721, -0.05160000840635173, -0.037176157641708965, -0.02574419680708342, 0.008411785616846546, 0.03798418337443638, -0.02658413754945694, -0.012150466297372464, -0.03380283348437231, -0.052287287966788186, -0.03304009609911819]"
$ws.Range("B7").Value = 0.01079919086965993
$ws.Range("B9").Value = 0.7883243645650251
